$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Duplicate the "Test2" block (rows 4-6) into a brand new "Test3" block
#    (rows 7-9), preserving formatting exactly via Copy (keeps style indices
#    identical to the source cells).
# ---------------------------------------------------------------------------
$ws.Range("A4:L6").Copy($ws.Range("A7:L9"))

# The source block (A4:L6) has some genuinely empty cells inside the
# rectangular copy area; Copy() materializes them as empty-but-present <c>
# elements in the destination. Remove those so row 7/8/9 only contain the
# cells that exist in the target layout.
$ws.Range("L7").ClearContents()
$ws.Range("A8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("A9:B9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("I9:K9").ClearContents()

# ---------------------------------------------------------------------------
# 2. Rename the new block's scenario markers from "Test2" to "Test3".
#    (Plain assignment - no leading apostrophe - keeps the fill-only style
#    used by A4/L6 instead of adding a quote-prefixed variant.)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Test3"
$ws.Range("L9").Value = "Test3"

# ---------------------------------------------------------------------------
# 3. Refresh the "dynamic" numeric id values. These are stored as text
#    (quote-prefixed shared strings), so a leading apostrophe is used to
#    force text instead of falling back to a numeric cell type.
#    Order matters for shared-string table layout: touch row 8 first,
#    right-to-left within the row, then mirror the same values into row 5
#    and row 2.
# ---------------------------------------------------------------------------
$ws.Range("K8").Value = "'749936"
$ws.Range("I8").Value = "'749936"
$ws.Range("H8").Value = "'749933"
$ws.Range("F8").Value = "'749933"
$ws.Range("E8").Value = "'749939"
$ws.Range("C8").Value = "'749939"

$ws.Range("K5").Value = "'749936"
$ws.Range("I5").Value = "'749936"
$ws.Range("H5").Value = "'749933"
$ws.Range("F5").Value = "'749933"
$ws.Range("E5").Value = "'749939"
$ws.Range("C5").Value = "'749939"

$ws.Range("E2").Value = "'749939"
$ws.Range("C2").Value = "'749939"

# ---------------------------------------------------------------------------
# 4. Row 2 gains three additional (blank, formatted) trailing cells so it
#    matches the look of rows 5/8, copied from the equivalent blank cells
#    in row 6 to keep the same style.
# ---------------------------------------------------------------------------
$ws.Range("F6:H6").Copy($ws.Range("F2:H2"))

# ---------------------------------------------------------------------------
# 5. New hyperlink on the new block's URL cell (B8), mirroring B2/B5.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B8"), "https://localhost:8080/") | Out-Null

# ---------------------------------------------------------------------------
# 6. Column width tweaks (col C and col F got narrower).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 12.28515625

# ---------------------------------------------------------------------------
# 7. Sheet dimension / active selection follow the grown data range.
# ---------------------------------------------------------------------------
$ws.Range("E12").Select()

Write-Host "Edit applied"
